# BIS-768: Fixed XLS export test files
# Adds a new "Unique" column (L) to the DATASET_TYPE export sheet, mirroring
# the existing "Multivalued" column (K): a bold header in row 4 and a
# TRUE/FALSE-styled value in rows 5-7.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Make room for the new column right after "Multivalued" (column K).
$ws.Columns("L").Insert()

# Reuse the "Multivalued" header/body formatting for the new "Unique" column.
$ws.Range("K4").Copy()
$ws.Range("L4").PasteSpecial(-4122)

$ws.Range("K5:K7").Copy()
$ws.Range("L5:L7").PasteSpecial(-4122)

# Fill in the new column's content. The data rows use the same literal
# "TRUE"/"FALSE" text convention as the "Multivalued" column, so a leading
# apostrophe keeps them as text instead of being auto-typed as booleans.
$ws.Range("L4").Value = "Unique"
$ws.Range("L5").Value = "'FALSE"
$ws.Range("L6").Value = "'FALSE"
$ws.Range("L7").Value = "'FALSE"

# Match the selection left behind by the edit.
$ws.Range("L4:L7").Select()
